$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(44308, 0, 4, 66.54466810846782),
    @(44309, 1, 5, 83.18083513558476),
    @(44310, 0, 4, 66.54466810846782),
    @(44311, 2, 6, 99.81700216270171),
    @(44312, 1, 7, 116.4531691898187)
)

$startRow = 234
$templateRow = 233

# Copy the formatting of the last existing row down onto the new rows
$endRow = $startRow + $data.Length - 1
$ws.Range("A$templateRow").Copy()
$ws.Range("A${startRow}:A${endRow}").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]

    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
}
